# Updates the "Avverkningsanmälningar" sheet:
#  - The "Förändrad" (column C) deadline date moves from 2026-02-13 (46066)
#    to 2026-02-17 (46070) for every data row (2-16).
#  - Row 14 ("A 58926-2025") becomes the new row 11, pushing the former
#    rows 11-13 down by one (11->12, 12->13, 13->14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every "Förändrad" date (column C, rows 2-16) forward by 4 days.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46070
}

# Capture the current (pre-reorder) values for rows 11-14 that are about to move.
# NOTE: use .Value2 for reads - .Value (getter) is unreliable in this runtime.
$rowA11 = $ws.Cells.Item(11, 1).Value2
$rowB11 = $ws.Cells.Item(11, 2).Value2
$rowG11 = $ws.Cells.Item(11, 7).Value2

$rowA12 = $ws.Cells.Item(12, 1).Value2
$rowB12 = $ws.Cells.Item(12, 2).Value2
$rowG12 = $ws.Cells.Item(12, 7).Value2

$rowA13 = $ws.Cells.Item(13, 1).Value2
$rowB13 = $ws.Cells.Item(13, 2).Value2
$rowG13 = $ws.Cells.Item(13, 7).Value2

$rowA14 = $ws.Cells.Item(14, 1).Value2
$rowB14 = $ws.Cells.Item(14, 2).Value2
$rowG14 = $ws.Cells.Item(14, 7).Value2

# Row 11 now holds what used to be row 14.
$ws.Cells.Item(11, 1).Value = $rowA14
$ws.Cells.Item(11, 2).Value = $rowB14
$ws.Cells.Item(11, 7).Value = $rowG14

# Row 12 now holds what used to be row 11.
$ws.Cells.Item(12, 1).Value = $rowA11
$ws.Cells.Item(12, 2).Value = $rowB11
$ws.Cells.Item(12, 7).Value = $rowG11

# Row 13 now holds what used to be row 12.
$ws.Cells.Item(13, 1).Value = $rowA12
$ws.Cells.Item(13, 2).Value = $rowB12
$ws.Cells.Item(13, 7).Value = $rowG12

# Row 14 now holds what used to be row 13.
$ws.Cells.Item(14, 1).Value = $rowA13
$ws.Cells.Item(14, 2).Value = $rowB13
$ws.Cells.Item(14, 7).Value = $rowG13
